$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) OOP paragraph: "... = 14.9 + 71 @ 10% = 22% + ...@15%"
#    - "14.9" loses its bold formatting
#    - "22%" gains bold formatting (note it is currently merged into the
#      run " = 22%"; setting Bold on just "22%" splits the run so the
#      leading " = " stays un-bold)
# ------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("14.9", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Bold = 0

$rng = $d.Content
$rng.Find.Execute("22%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Bold = 1

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark out of the OOP paragraph (it used to sit
#    between " = 22%" and " + ...@15%") into the following, otherwise
#    empty, ListParagraph paragraph.
# ------------------------------------------------------------------

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$targetPara = $d.Paragraphs(3).Range
$d.Bookmarks.Add("_GoBack", $targetPara)

# ------------------------------------------------------------------
# 3) "Software eng 2" paragraph: wrap the two terms in parentheses and
#    append the new " + (83 @15% = 12.45)" term, ending with the bold
#    subtotal "17.625".
# ------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute(" 2: 69 @ 7.5% = ", $true, $false, $false, $false, $false, $true, 1, $false, " 2: (69 @ 7.5% = ", 2)

$rng = $d.Content
$rng.Find.Execute(" + …@15% + …@ 7.5%", $true, $false, $false, $false, $false, $true, 1, $false, ") + (83 @15% = 12.45) = 17.625 + …@ 7.5%", 2)

# Bold the ")" that immediately follows "5.175"
$rng = $d.Content
$rng.Find.Execute("5.175)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$closeParen = $d.Range($rng.End - 1, $rng.End)
$closeParen.Font.Bold = 1

# Bold the new subtotal "17.625"
$rng = $d.Content
$rng.Find.Execute("17.625", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Bold = 1

Write-Host "done"
